# Bug fix in Kbar now-casting calculation.
# Updates columns C (E_K_fc_fin) and D (E_K_fc_capital_fin) for rows 2-50
# on Sheet1. Column E (E_K_fc_other_fin = C - D) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$data = @(
    @(2, 81.09106943735266, 10.13996617324147),
    @(3, 98.4204690087752, 12.02630076831704),
    @(4, 119.4429166398555, 13.07043180840359),
    @(5, 27.89256489766104, 2.976109206729212),
    @(6, 25.87783081391243, 1.636169826518793),
    @(7, 181.8035109186626, 64.6037853725084),
    @(8, 91.80975688474867, 17.66595338990206),
    @(9, 20.52188610557545, 4.209222070133255),
    @(10, 62.86138988194384, 5.276951174403877),
    @(11, 53.66290083119358, 8.614811797120455),
    @(12, 88.37959858545823, 10.16349369387679),
    @(13, 109.6919803657956, 13.13572435349594),
    @(14, 51.74872726013162, 6.732329609910158),
    @(15, 44.88764127486557, 5.866722049813027),
    @(16, 114.9389027574633, 16.12959376804319),
    @(17, 66.47039923590675, 8.036246018373033),
    @(18, 64.98130380093853, 7.749892194827897),
    @(19, 115.5336202580731, 6.977014376054265),
    @(20, 36.64975370606664, 4.13412300085932),
    @(21, 39.70211150743625, 5.708687582773737),
    @(22, 10.41642217933556, 1.05090379488537),
    @(23, 73.59300222877866, 8.462344458858782),
    @(24, 6.272130783095512, 0.5365625781674238),
    @(25, 58.66873979331869, 7.362495051132077),
    @(26, 88.37801185544566, 34.08576114963231),
    @(27, 69.34879879626797, 13.0062631300152),
    @(28, 55.71036424334256, 4.517764755528603),
    @(29, 190.2666917046467, 18.42645392363421),
    @(30, 63.22384352504145, 7.497682038862504),
    @(31, 69.8647025646125, 7.061893287011926),
    @(32, 17.49108955688107, 1.895466618349366),
    @(33, 99.52386531570954, 10.46408872112913),
    @(34, 145.7929087768983, 20.06065309201589),
    @(35, 43.96739964255068, 4.780635356628347),
    @(36, 45.53085768534857, 5.281583074836819),
    @(37, 19.65123610670828, 2.384439373554106),
    @(38, 45.10051721391826, 3.857730944760764),
    @(39, 87.07867503560531, 11.10764327160515),
    @(40, 48.08977013335183, 5.530937271576509),
    @(41, 48.00476824953135, 6.591683050160187),
    @(42, 32.09902297667693, 3.459031476921247),
    @(43, 67.04340374934132, 4.063472752867128),
    @(44, 105.3052184104982, 12.22429155675636),
    @(45, 12.84495941074148, 1.765448946907162),
    @(46, 27.85905107255654, 1.625571266440792),
    @(47, 4.607196448654674, 0.5243634985325216),
    @(48, 26.62102564419216, 2.691638575206132),
    @(49, 33.38746513947333, 3.91018825707657),
    @(50, 25.44325883016047, 1.986349047012836)
)

foreach ($item in $data) {
    $row = $item[0]
    $cValue = $item[1]
    $dValue = $item[2]
    $ws.Cells.Item($row, 3).Value = $cValue
    $ws.Cells.Item($row, 4).Value = $dValue
}
